# SALES improved contagents 0.8 at 15/10/2022
# Refresh the tyre sales table (sheet "Holidays 2019") columns E:J
# - rows 2-19/20-39: updated Date_of_sales / Contragent (and a few Sales value corrections)
# - rows 40-56: newly appended sales records
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Tyre Size, Model, Param, Sales value, Date_of_sales (serial), Contragent
$data = @(
  ,("315/80R22.5", "BEL-158M", "камневыт, груз, сер", 259, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-158M", "камневыт, груз, трп", 285, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-278", "груз, сер", 281, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-278", "груз, трп", 287, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-268", "груз, сер", 283, 44814, "БНХ РОС")
  ,("315/80R22.5", "BEL-268", "груз, трп", 287, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-398", "груз, сер", 265, 44814, "нет данных")
  ,("315/80R22.5", "BEL-326", "груз, сер", 269, 44814, "нет данных")
  ,("315/80R22.5", "BEL-326", "груз, трп", 27, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-498", "156L, груз, сер", 247, 44814, "нет данных")
  ,("315/80R22.5", "BEL-518", "груз, сер", 255, 44814, "HHHDFD")
  ,("12.00R20", "ИД-304М", "16, груз, сер", 1555, 44814, "нет данных")
  ,("12.00R20", "ИД-304М", "18, груз, сер", 285, 44814, "HHHDFD")
  ,("12.00R20", "ИД-304М", "16, груз, трп", 290, 44814, "HHHDFD")
  ,("12.00R20", "ИД-304М", "18, груз, трп", 296, 44814, "HHHDFD")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 4565, 44814, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 4565, 44814, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 283, 44814, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 283, 44814, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, трп", 291, 44814, "HHHDFD")
  ,("195/65R15", "BEL-337", "б/к, сер, легк", 291, 44814, "HHHDFD")
  ,("315/80R22.5", "BEL-158M", "камневыт, груз, сер", 259, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-158M", "камневыт, груз, трп", 285, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-278", "груз, сер", 281, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-278", "груз, трп", 287, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-268", "груз, сер", 283, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-268", "груз, трп", 287, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-398", "груз, сер", 265, 44752, "нет данных")
  ,("315/80R22.5", "BEL-326", "груз, сер", 269, 44752, "нет данных")
  ,("315/80R22.5", "BEL-326", "груз, трп", 27, 44752, "БНХ РОС")
  ,("315/80R22.5", "BEL-498", "156L, груз, сер", 247, 44752, "нет данных")
  ,("315/80R22.5", "BEL-518", "груз, сер", 255, 44752, "БНХ РОС")
  ,("12.00R20", "ИД-304М", "16, груз, сер", 1555, 44752, "нет данных")
  ,("12.00R20", "ИД-304М", "18, груз, сер", 285, 44752, "БНХ РОС")
  ,("12.00R20", "ИД-304М", "16, груз, трп", 290, 44752, "БНХ РОС")
  ,("12.00R20", "ИД-304М", "18, груз, трп", 296, 44752, "БНХ РОС")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 4565, 44752, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 4565, 44752, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 283, 44752, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, сер", 283, 44752, "нет данных")
  ,("12.00R20", "БИ-368М", "18, груз, трп", 291, 44752, "HHHDFD")
  ,("195/65R15", "BEL-337", "б/к, сер, легк", 291, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-158M", "камневыт, груз, сер", 259, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-158M", "камневыт, груз, трп", 285, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-278", "груз, сер", 281, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-278", "груз, трп", 287, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-268", "груз, сер", 283, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-268", "груз, трп", 287, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-326", "груз, трп", 27, 44752, "HHHDFD")
  ,("315/80R22.5", "BEL-518", "груз, сер", 255, 44752, "HHHDFD")
  ,("12.00R20", "ИД-304М", "18, груз, сер", 285, 44752, "HHHDFD")
  ,("12.00R20", "ИД-304М", "16, груз, трп", 290, 44752, "БНХ УКР")
  ,("12.00R20", "ИД-304М", "18, груз, трп", 296, 44752, "БНХ УКР")
  ,("12.00R20", "БИ-368М", "18, груз, трп", 291, 44752, "БНХ УКР")
  ,("195/65R15", "BEL-337", "б/к, сер, легк", 291, 44752, "БНХ УКР")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2   # data starts at row 2 (row 1 is the header)
  $row = $data[$i]
  $ws.Cells.Item($r, 5).Value  = $row[0]   # E Tyre Size
  $ws.Cells.Item($r, 6).Value  = $row[1]   # F Model
  $ws.Cells.Item($r, 7).Value  = $row[2]   # G Param
  $ws.Cells.Item($r, 8).Value  = $row[3]   # H Sales value
  $ws.Cells.Item($r, 9).Value  = $row[4]   # I Date_of_sales
  $ws.Cells.Item($r, 10).Value = $row[5]   # J Contragent
}
